# Add the new "Step1 / Step 2 / definition of equality" content to the
# workflow section, right after the embedded Visio object paragraph and
# before the trailing bookmark paragraph.

$d = $word.ActiveDocument

# Paragraph that holds the embedded Visio object (w:object) - new content
# is inserted right after it.
$anchorIndex = 3

function Insert-ParaWithRuns {
    param([int]$afterIndex, [string]$styleName, [string[]]$runs)

    # Create a brand new paragraph right after $afterIndex.
    $afterPara = $d.Paragraphs.Item($afterIndex)
    $afterPara.Range.InsertParagraphAfter()
    $newIndex = $afterIndex + 1
    $newPara = $d.Paragraphs.Item($newIndex)

    if ($styleName) {
        $newPara.Style = $styleName
    } else {
        $newPara.Style = "Normal"
    }

    # Put the first run of text directly into the new paragraph.
    $newPara.Range.Text = $runs[0]

    # Every additional run is built by inserting ANOTHER paragraph after
    # the current one, filling it with that run's text, and then removing
    # the paragraph mark that separates the two - this joins them back
    # into a single paragraph while keeping the runs distinct (the engine
    # merges same-format text typed/inserted into one run, but stitching
    # two paragraphs together preserves the run boundary).
    for ($i = 1; $i -lt $runs.Length; $i++) {
        $curPara = $d.Paragraphs.Item($newIndex)
        $curPara.Range.InsertParagraphAfter()
        $extraPara = $d.Paragraphs.Item($newIndex + 1)
        $extraPara.Range.Text = $runs[$i]

        $joinStart = $curPara.Range.End - 1
        $joinRange = $d.Range($joinStart, $joinStart + 1)
        $joinRange.Delete()
    }

    return $newIndex
}

$idx = $anchorIndex
$idx = Insert-ParaWithRuns $idx "Heading 1" @("Step1: Validation of the triangle")
$idx = Insert-ParaWithRuns $idx $null @("Three parameters are used to create a triangle. They have to meet some requirements before constructing the Triangle object.")
$idx = Insert-ParaWithRuns $idx $null @("Firstly, the lengths of three laterals must be greater than 0.")
$idx = Insert-ParaWithRuns $idx $null @("Secondly, the sum of any two laterals must be larger than the third one. Otherwise, it is not possible to form a triangle.")
$idx = Insert-ParaWithRuns $idx $null @("If the two criteria are not met, the system throw exceptions. It is better that the system warn the user and ask the user to re-input. However, for the time being, I would just throw the exception. It is to be enhanced in the future.")

$idx = Insert-ParaWithRuns $idx "Heading 1" @("Step 2: Decide the triangle type")
$idx = Insert-ParaWithRuns $idx $null @("There three types of triangle. The easiest one is equal-lateral. The condition for this is simple. If the first lateral equals the second and the second equals the third, the triangle is considered as equal-lateral.")

$idx = Insert-ParaWithRuns $idx $null @(
    "It is easier to determine whether a triangle is ",
    "scalene ",
    "than isosceles. Thus, the next sub-step is to check whether the triangle is ",
    "scalene",
    ". This needs three pairwise comparisons of the laterals, if all are different, the triangle is ",
    "scalene",
    ". Otherwise, it is ",
    "isosceles",
    "."
)

$idx = Insert-ParaWithRuns $idx "Heading 1" @("The definition of equality")

# The very last new paragraph does not get its own paragraph mark: its
# text lands inside the pre-existing trailing paragraph (the one holding
# the _GoBack bookmark), right before the bookmark.
$bookmarkPara = $d.Paragraphs.Item($idx + 1)
$bookmarkPara.Range.InsertBefore("In a real project, equality usefully is not 100% equality. Some error must be allowed. The error can be measure by length or percentage. I prefer to measure it by percentage. Percentage makes more sense. However, matters like that needs discussion with the customer.")

Write-Output "done"
